# Weekly update to the "Vega Monumental Concepción - Naranja" sheet.
# Two new price records (Naranja "Lane Late", Primera/Segunda) are added
# for the most recent week; all the older records shift down two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 121 - this pushes the
# existing rows 121:146 down to 123:148 (and extends the used range /
# dimension to A1:T148 automatically), exactly like a user inserting rows
# in the UI before typing in the new week's data.
$ws.Rows("121:122").Insert()

# --- New row 121: Naranja / Lane Late / Primera ---
$ws.Range("A121").Value = 11
$ws.Range("B121").Value = "Vega Monumental Concepción"
$ws.Range("C121").Value = "Bíobío"
$ws.Range("D121").Value = 44474
$ws.Range("E121").Value = 8
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100102
$ws.Range("H121").Value = "Cítricos"
$ws.Range("I121").Value = 100102005
$ws.Range("J121").Value = "Naranja"
$ws.Range("K121").Value = "Lane Late"
$ws.Range("L121").Value = "Primera"
$ws.Range("M121").Value = 100
$ws.Range("N121").Value = 7000
$ws.Range("O121").Value = 7500
$ws.Range("P121").Value = 7250
$ws.Range("Q121").Value = "$/caja 15 kilos empedrada"
$ws.Range("R121").Value = "Región de O'Higgins"
$ws.Range("S121").Value = 483
$ws.Range("T121").Value = 15

# --- New row 122: Naranja / Lane Late / Segunda ---
$ws.Range("A122").Value = 11
$ws.Range("B122").Value = "Vega Monumental Concepción"
$ws.Range("C122").Value = "Bíobío"
$ws.Range("D122").Value = 44474
$ws.Range("E122").Value = 8
$ws.Range("F122").Value = "Fruta"
$ws.Range("G122").Value = 100102
$ws.Range("H122").Value = "Cítricos"
$ws.Range("I122").Value = 100102005
$ws.Range("J122").Value = "Naranja"
$ws.Range("K122").Value = "Lane Late"
$ws.Range("L122").Value = "Segunda"
$ws.Range("M122").Value = 50
$ws.Range("N122").Value = 6500
$ws.Range("O122").Value = 6500
$ws.Range("P122").Value = 6500
$ws.Range("Q122").Value = "$/caja 15 kilos empedrada"
$ws.Range("R122").Value = "Región de O'Higgins"
$ws.Range("S122").Value = 433
$ws.Range("T122").Value = 15
